$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Add new header cell H1 with the "Numéro Siret" label
$ws.Range("H1").Value = "Numéro Siret"

# Update the selection to match the diff
$ws.Activate()
$ws.Range("G8").Select()
